$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ducks")

# Fix typo: Purchase_State "Nv" -> "NV" for row 109
$ws.Range("F109").Value = "NV"

# Add new "QQ" Purchase_State values for rows missing Purchase_State (column F)
$ws.Range("F14").Value = "QQ"
$ws.Range("F15").Value = "QQ"
$ws.Range("F23").Value = "QQ"
$ws.Range("F24").Value = "QQ"
$ws.Range("F25").Value = "QQ"
$ws.Range("F26").Value = "QQ"
$ws.Range("F27").Value = "QQ"
$ws.Range("F29").Value = "QQ"
$ws.Range("F72").Value = "QQ"
$ws.Range("F73").Value = "QQ"
$ws.Range("F74").Value = "QQ"
$ws.Range("F75").Value = "QQ"
$ws.Range("F76").Value = "QQ"
$ws.Range("F77").Value = "QQ"
$ws.Range("F78").Value = "QQ"
$ws.Range("F79").Value = "QQ"

# Fill in "About Me" ranking numbers (column L) for rows previously blank
$ws.Range("L41").Value = 1
$ws.Range("L44").Value = 2
$ws.Range("L47").Value = 3
$ws.Range("L48").Value = 4
$ws.Range("L49").Value = 5
$ws.Range("L51").Value = 6
$ws.Range("L53").Value = 7
$ws.Range("L54").Value = 8
$ws.Range("L63").Value = 9
$ws.Range("L64").Value = 10
$ws.Range("L65").Value = 11
$ws.Range("L69").Value = 12
$ws.Range("L70").Value = 13
$ws.Range("L72").Value = 14
$ws.Range("L73").Value = 15
$ws.Range("L74").Value = 16
$ws.Range("L75").Value = 17
$ws.Range("L76").Value = 18
$ws.Range("L77").Value = 19
$ws.Range("L78").Value = 20
$ws.Range("L79").Value = 21
$ws.Range("L80").Value = 22
$ws.Range("L81").Value = 23
$ws.Range("L82").Value = 24
$ws.Range("L83").Value = 25
$ws.Range("L84").Value = 26
$ws.Range("L85").Value = 27
$ws.Range("L86").Value = 28
$ws.Range("L89").Value = 29
$ws.Range("L91").Value = 30
$ws.Range("L92").Value = 31
$ws.Range("L93").Value = 32
$ws.Range("L94").Value = 33
$ws.Range("L97").Value = 34
$ws.Range("L98").Value = 35
$ws.Range("L99").Value = 36
$ws.Range("L100").Value = 37
$ws.Range("L101").Value = 38
$ws.Range("L102").Value = 39
$ws.Range("L103").Value = 40
$ws.Range("L104").Value = 41
$ws.Range("L105").Value = 42
$ws.Range("L106").Value = 43
$ws.Range("L107").Value = 44
$ws.Range("L108").Value = 45
$ws.Range("L109").Value = 46

# Fill in Name/ranking numbers (column B) for rows previously blank
$ws.Range("B48").Value = 21
$ws.Range("B53").Value = 22
$ws.Range("B59").Value = 1
$ws.Range("B60").Value = 2
$ws.Range("B66").Value = 3
$ws.Range("B67").Value = 4
$ws.Range("B80").Value = 5
$ws.Range("B82").Value = 6
$ws.Range("B83").Value = 7
$ws.Range("B86").Value = 8
$ws.Range("B87").Value = 9
$ws.Range("B97").Value = 20
$ws.Range("B100").Value = 10
$ws.Range("B101").Value = 11
$ws.Range("B102").Value = 12
$ws.Range("B103").Value = 13
$ws.Range("B104").Value = 14
$ws.Range("B105").Value = 15
$ws.Range("B106").Value = 16
$ws.Range("B107").Value = 17
$ws.Range("B108").Value = 18
$ws.Range("B109").Value = 19

# Update the active view/selection to match final state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 92
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("F109").Select()
